$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 4 rows (Resolving-Mac as sender) which no longer appear in the updated data
$ws.Range("A14:T17").Delete() | Out-Null

# Update values recomputed with the new TPM-based data (columns E-T) for the remaining rows

# Row 2
$ws.Range("G2").Value = 1.008000333333333
$ws.Range("H2").Value = 3.024001
$ws.Range("I2").Value = 0.3525296793986107
$ws.Range("J2").Value = 0.3525296793986107
$ws.Range("M2").Value = 19.163974
$ws.Range("N2").Value = 57.491922
$ws.Range("O2").Value = 0.6845732287637933
$ws.Range("P2").Value = 0.6845732287637933
$ws.Range("Q2").Value = 19.31729217999133
$ws.Range("R2").Value = 173.855629619922
$ws.Range("S2").Value = 0.2413323808609719
$ws.Range("T2").Value = 0.2413323808609718

# Row 3
$ws.Range("G3").Value = 1.008000333333333
$ws.Range("H3").Value = 3.024001
$ws.Range("I3").Value = 0.3525296793986107
$ws.Range("J3").Value = 0.3525296793986107
$ws.Range("O3").Value = 0.02733363438148322
$ws.Range("P3").Value = 0.02733363438148323
$ws.Range("Q3").Value = 0.7713006870596667
$ws.Range("R3").Value = 6.941706183537
$ws.Range("S3").Value = 0.009635917365303124
$ws.Range("T3").Value = 0.009635917365303124

# Row 4
$ws.Range("G4").Value = 1.008000333333333
$ws.Range("H4").Value = 3.024001
$ws.Range("I4").Value = 0.3525296793986107
$ws.Range("J4").Value = 0.3525296793986107
$ws.Range("M4").Value = 7.880893333333333
$ws.Range("N4").Value = 23.64268
$ws.Range("O4").Value = 0.281520346184098
$ws.Range("P4").Value = 0.281520346184098
$ws.Range("Q4").Value = 7.943943106964444
$ws.Range("R4").Value = 71.49548796268
$ws.Range("S4").Value = 0.09924427738446598
$ws.Range("T4").Value = 0.09924427738446595

# Row 5
$ws.Range("G5").Value = 1.008000333333333
$ws.Range("H5").Value = 3.024001
$ws.Range("I5").Value = 0.3525296793986107
$ws.Range("J5").Value = 0.3525296793986107
$ws.Range("M5").Value = 0.183999
$ws.Range("N5").Value = 0.551997
$ws.Range("O5").Value = 0.006572790670625477
$ws.Range("P5").Value = 0.006572790670625476
$ws.Range("Q5").Value = 0.185471053333
$ws.Range("R5").Value = 1.669239479997
$ws.Range("S5").Value = 0.002317103787869779
$ws.Range("T5").Value = 0.002317103787869779

# Row 6
$ws.Range("I6").Value = 0.5377259091975243
$ws.Range("J6").Value = 0.5377259091975243
$ws.Range("M6").Value = 19.163974
$ws.Range("N6").Value = 57.491922
$ws.Range("O6").Value = 0.6845732287637933
$ws.Range("P6").Value = 0.6845732287637933
$ws.Range("Q6").Value = 29.46534464400333
$ws.Range("R6").Value = 265.18810179603
$ws.Range("S6").Value = 0.3681127618492955
$ws.Range("T6").Value = 0.3681127618492955

# Row 7
$ws.Range("I7").Value = 0.5377259091975243
$ws.Range("J7").Value = 0.5377259091975243
$ws.Range("O7").Value = 0.02733363438148322
$ws.Range("P7").Value = 0.02733363438148323
$ws.Range("S7").Value = 0.01469800339945578
$ws.Range("T7").Value = 0.01469800339945578

# Row 8
$ws.Range("I8").Value = 0.5377259091975243
$ws.Range("J8").Value = 0.5377259091975243
$ws.Range("M8").Value = 7.880893333333333
$ws.Range("N8").Value = 23.64268
$ws.Range("O8").Value = 0.281520346184098
$ws.Range("P8").Value = 0.281520346184098
$ws.Range("Q8").Value = 12.11717560091111
$ws.Range("R8").Value = 109.0545804082
$ws.Range("S8").Value = 0.1513807841094459
$ws.Range("T8").Value = 0.1513807841094459

# Row 9
$ws.Range("I9").Value = 0.5377259091975243
$ws.Range("J9").Value = 0.5377259091975243
$ws.Range("M9").Value = 0.183999
$ws.Range("N9").Value = 0.551997
$ws.Range("O9").Value = 0.006572790670625477
$ws.Range("P9").Value = 0.006572790670625476
$ws.Range("Q9").Value = 0.282905515795
$ws.Range("R9").Value = 2.546149642155
$ws.Range("S9").Value = 0.00353435983932709
$ws.Range("T9").Value = 0.00353435983932709

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.313796
$ws.Range("H10").Value = 0.941388
$ws.Range("I10").Value = 0.1097444114038651
$ws.Range("J10").Value = 0.1097444114038651
$ws.Range("M10").Value = 19.163974
$ws.Range("N10").Value = 57.491922
$ws.Range("O10").Value = 0.6845732287637933
$ws.Range("P10").Value = 0.6845732287637933
$ws.Range("Q10").Value = 6.013578385304
$ws.Range("R10").Value = 54.122205467736
$ws.Range("S10").Value = 0.07512808605352597
$ws.Range("T10").Value = 0.07512808605352596

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.313796
$ws.Range("H11").Value = 0.941388
$ws.Range("I11").Value = 0.1097444114038651
$ws.Range("J11").Value = 0.1097444114038651
$ws.Range("O11").Value = 0.02733363438148322
$ws.Range("P11").Value = 0.02733363438148323
$ws.Range("Q11").Value = 0.240110109484
$ws.Range("R11").Value = 2.160990985356
$ws.Range("S11").Value = 0.002999713616724325
$ws.Range("T11").Value = 0.002999713616724325

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.313796
$ws.Range("H12").Value = 0.941388
$ws.Range("I12").Value = 0.1097444114038651
$ws.Range("J12").Value = 0.1097444114038651
$ws.Range("M12").Value = 7.880893333333333
$ws.Range("N12").Value = 23.64268
$ws.Range("O12").Value = 0.281520346184098
$ws.Range("P12").Value = 0.281520346184098
$ws.Range("Q12").Value = 2.472992804426667
$ws.Range("R12").Value = 22.25693523984
$ws.Range("S12").Value = 0.03089528469018616
$ws.Range("T12").Value = 0.03089528469018616

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.313796
$ws.Range("H13").Value = 0.941388
$ws.Range("I13").Value = 0.1097444114038651
$ws.Range("J13").Value = 0.1097444114038651
$ws.Range("M13").Value = 0.183999
$ws.Range("N13").Value = 0.551997
$ws.Range("O13").Value = 0.006572790670625477
$ws.Range("P13").Value = 0.006572790670625476
$ws.Range("Q13").Value = 0.057738150204
$ws.Range("R13").Value = 0.519643351836
$ws.Range("S13").Value = 0.0007213270434286086
$ws.Range("T13").Value = 0.0007213270434286084
